# Pewlett-Hackard DBD schema: mark Titles.title / Titles.from_date and
# Managers.emp_no / Managers.from_date as part of the primary key.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Managers table (rows 36-41) ---
$ws.Range("A39").Value = "emp_no fk - Employees.emp_no pk "
$ws.Range("A40").Value = "from_date date "

# --- Titles table (rows 22-27) ---
$ws.Range("A25").Value = "title varchar pk"
$ws.Range("A26").Value = "from_date date pk"

# Update the view state to match where the author left the cursor/scroll
# position after making the edit.
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 1
$ws.Range("C26").Select()
